$wb = $excel.ActiveWorkbook

# --- Sheet "Trends Status": Insufficient Data row updated from 299 -> 300 ---
$wsTrends = $wb.Worksheets.Item("Trends Status")
$wsTrends.Range("B8").Value = 300
$wsTrends.Range("C8").Value = 300

# --- Sheet "Priority Status": counts updated ---
$wsPriority = $wb.Worksheets.Item("Priority Status")
$wsPriority.Range("B2").Value = 103
$wsPriority.Range("B3").Value = 286
$wsPriority.Range("B4").Value = 554

# --- Sheet "Species qualification": label + count updated ---
$wsQual = $wb.Worksheets.Item("Species qualification")
$wsQual.Range("A2").Value = "SoIB Assessment"
$wsQual.Range("B2").Value = 300

# --- Sheet "High Priority break-up" handling ---
# The original sheet (sheetId 5) keeps its identity but is renamed to
# "Interannual update - High Pri" and receives new break-up data.
# A brand new sheet (sheetId 6), "Major update - High Priority ", is
# inserted right after it, preserving the *old* break-up data untouched.

$wsOld = $wb.Worksheets.Item("High Priority break-up")

# Add the new sheet right after the existing one, before renaming/rewriting it,
# so the new sheet inherits the old content first.
$wsNew = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsOld)
$wsNew.Name = "Major update - High Priority "

# Copy over the old sheet's contents into the new sheet verbatim.
$wsNew.Range("A1").Value = "Break-up"
$wsNew.Range("B1").Value = "High Species (no.)"
$wsNew.Range("C1").Value = "High Species (perc.)"
$wsNew.Range("D1").Value = "New High Species (no.)"
$wsNew.Range("E1").Value = "New High Species (perc.)"

# Reuse the existing bold+centered header formatting (style already present
# in the workbook) instead of re-building it step by step, which would leave
# behind an orphaned intermediate style.
$wsOld.Range("A1:E1").Copy()
$wsNew.Range("A1:E1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wsNew.Range("A2").Value = "IUCN"
$wsNew.Range("B2").Value = 4
$wsNew.Range("C2").Value = 100
$wsNew.Range("D2").Value = 4
$wsNew.Range("E2").Value = 100

# Now rename the original sheet and overwrite it with the new break-up data.
$wsOld.Name = "Interannual update - High Pri"

$wsOld.Range("A2").Value = "Trend New"
$wsOld.Range("B2").Value = 86
$wsOld.Range("C2").Value = 83.5
$wsOld.Range("D2").Value = 86
$wsOld.Range("E2").Value = 86

$wsOld.Range("A3").Value = "IUCN"
$wsOld.Range("B3").Value = 17
$wsOld.Range("C3").Value = 16.5
$wsOld.Range("D3").Value = 14
$wsOld.Range("E3").Value = 14
